$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (row 1): Wins, Losses, Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting from an existing header cell (A1) so the
# new header cells pick up the same bold/centered/bordered style.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the season record (Wins/Losses/Ties) for every data row.
for ($r = 2; $r -le 46; $r++) {
    $ws.Range("AD$r").Value = 76
    $ws.Range("AE$r").Value = 86
    $ws.Range("AF$r").Value = 0
}
